$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value2 = 1135658.314116
$ws.Cells.Item(3, 4).Value2 = 9454098.163177
$ws.Cells.Item(4, 4).Value2 = 3655125.363434
$ws.Cells.Item(5, 4).Value2 = 1232183.189674
$ws.Cells.Item(6, 4).Value2 = 4838313.764562
$ws.Cells.Item(7, 4).Value2 = 1285552.33784
$ws.Cells.Item(8, 4).Value2 = 368679.286824
$ws.Cells.Item(9, 4).Value2 = 608831.68731
$ws.Cells.Item(10, 4).Value2 = 964906.804544
$ws.Cells.Item(11, 4).Value2 = 1095851.850221
$ws.Cells.Item(12, 4).Value2 = 840683.441916
$ws.Cells.Item(13, 4).Value2 = 2164284.563868
$ws.Cells.Item(14, 4).Value2 = 1097109.467157
$ws.Cells.Item(15, 4).Value2 = 469841.620535
$ws.Cells.Item(16, 4).Value2 = 5979774.616828
$ws.Cells.Item(17, 4).Value2 = 1077970.096821
$ws.Cells.Item(18, 4).Value2 = 2482136.912266
$ws.Cells.Item(19, 4).Value2 = 1583402.958161
$ws.Cells.Item(20, 4).Value2 = 330579.926203
$ws.Cells.Item(21, 4).Value2 = 558258.659549
$ws.Cells.Item(22, 4).Value2 = 74691253.674562
$ws.Cells.Item(23, 4).Value2 = 2497618.090538
$ws.Cells.Item(24, 4).Value2 = 990081.85941
$ws.Cells.Item(25, 4).Value2 = 3620612.016261
$ws.Cells.Item(26, 4).Value2 = 4319821.975134
$ws.Cells.Item(27, 4).Value2 = 2705992.625732
$ws.Cells.Item(28, 4).Value2 = 2360529.146302
$ws.Cells.Item(29, 4).Value2 = 72179.732786
$ws.Cells.Item(30, 4).Value2 = 1307378.460047
$ws.Cells.Item(31, 4).Value2 = 1010724.059407
$ws.Cells.Item(32, 4).Value2 = 1697662.641165
$ws.Cells.Item(33, 4).Value2 = 1624915.22003
